$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formatting (borders/fills) from the existing pattern blocks onto the new ranges ---
$ws.Range("A3:G9").Copy()
$ws.Range("A24:G30").PasteSpecial(-4122)
$ws.Range("A12:G18").Copy()
$ws.Range("A33:G39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New header row 22 (text labels) ---
$ws.Range("A22").Value = "C"
$ws.Range("B22").Value = "S"

# Row 23
$ws.Cells.Item(23,1).Value = 8
$ws.Cells.Item(23,2).Value = 1
$ws.Cells.Item(23,3).Value = 10
$ws.Cells.Item(23,4).Value = 2
$ws.Cells.Item(23,5).Value = 8
$ws.Cells.Item(23,6).Value = 6
$ws.Cells.Item(23,7).Value = 4

# Row 24
$ws.Cells.Item(24,2).Value = 0
$ws.Cells.Item(24,3).Value = 1
$ws.Cells.Item(24,4).Value = 2
$ws.Cells.Item(24,5).Value = 3
$ws.Cells.Item(24,6).Value = 4
$ws.Cells.Item(24,7).Value = 5

# Row 25
$ws.Cells.Item(25,1).Value = 0

# Row 26
$ws.Cells.Item(26,1).Value = 1
$ws.Cells.Item(26,2).Formula = "=MAX(B27+B35,C27)"

# Row 27
$ws.Cells.Item(27,1).Value = 2
$ws.Cells.Item(27,2).Value = 0
$ws.Cells.Item(27,3).Formula = "=MAX(D28,`$B`$28+C36)"

# Row 28
$ws.Cells.Item(28,1).Value = 3
$ws.Cells.Item(28,2).Value = 0
$ws.Cells.Item(28,3).Formula = "=MAX(C37+`$B`$29,D29)"
$ws.Cells.Item(28,4).Formula = "=MAX(D37+`$B`$29,E29)"

# Row 29
$ws.Cells.Item(29,1).Value = 4
$ws.Cells.Item(29,2).Value = 0
$ws.Cells.Item(29,3).Formula = "=MAX(D30,C38+`$B`$30)"
$ws.Cells.Item(29,4).Formula = "=MAX(E30,D38+`$B`$30)"
$ws.Cells.Item(29,5).Formula = "=MAX(F30,E38+`$B`$30)"

# Row 30
$ws.Cells.Item(30,1).Value = 5
$ws.Cells.Item(30,2).Formula = "=MAX(B39,0)"
$ws.Cells.Item(30,3).Formula = "=C39"
$ws.Cells.Item(30,4).Formula = "=D39"
$ws.Cells.Item(30,5).Formula = "=E39"
$ws.Cells.Item(30,6).Formula = "=F39"

# Row 33
$ws.Cells.Item(33,2).Value = 0
$ws.Cells.Item(33,3).Value = 1
$ws.Cells.Item(33,4).Value = 2
$ws.Cells.Item(33,5).Value = 3
$ws.Cells.Item(33,6).Value = 4
$ws.Cells.Item(33,7).Value = 5

# Row 34
$ws.Cells.Item(34,1).Value = 0

# Row 35
$ws.Cells.Item(35,1).Value = 1
$ws.Cells.Item(35,2).Formula = "=C35-`$A`$23"
$ws.Cells.Item(35,3).Value = 10
$ws.Cells.Item(35,4).Formula = "=C35*(1-`$B`$23)"
$ws.Cells.Item(35,5).Formula = "=D35*(1-`$B`$23)"
$ws.Cells.Item(35,6).Formula = "=E35*(1-`$B`$23)"
$ws.Cells.Item(35,7).Formula = "=F35*(1-`$B`$23)"

# Row 36
$ws.Cells.Item(36,1).Value = 2
$ws.Cells.Item(36,2).Formula = "=D36-A23"
$ws.Cells.Item(36,3).Formula = "=B36+D35"
$ws.Cells.Item(36,4).Value = 2
$ws.Cells.Item(36,5).Formula = "=D36*(1-`$B`$23)"
$ws.Cells.Item(36,6).Formula = "=E36*(1-`$B`$23)"
$ws.Cells.Item(36,7).Formula = "=F36*(1-`$B`$23)"

# Row 37
$ws.Cells.Item(37,1).Value = 3
$ws.Cells.Item(37,2).Formula = "=E37-A23"
$ws.Cells.Item(37,3).Formula = "=B37+E36"
$ws.Cells.Item(37,4).Formula = "=C37+E35"
$ws.Cells.Item(37,5).Value = 8
$ws.Cells.Item(37,6).Formula = "=E37*(1-`$B`$23)"
$ws.Cells.Item(37,7).Formula = "=F37*(1-`$B`$23)"

# Row 38
$ws.Cells.Item(38,1).Value = 4
$ws.Cells.Item(38,2).Formula = "=F38-A23"
$ws.Cells.Item(38,3).Formula = "=B38+F37"
$ws.Cells.Item(38,4).Formula = "=C38+F36"
$ws.Cells.Item(38,5).Formula = "=D38+F35"
$ws.Cells.Item(38,6).Value = 6
$ws.Cells.Item(38,7).Formula = "=F38*(1-`$B`$23)"

# Row 39
$ws.Cells.Item(39,1).Value = 5
$ws.Cells.Item(39,2).Formula = "=G39-`$A`$23"
$ws.Cells.Item(39,3).Formula = "=B39+G38"
$ws.Cells.Item(39,4).Formula = "=C39+G37"
$ws.Cells.Item(39,5).Formula = "=D39+G36"
$ws.Cells.Item(39,6).Formula = "=E39+G35"
$ws.Cells.Item(39,7).Value = 4
# --- Final selection state ---
$ws.Range("B30").Select()
